$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2110091743119266
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.01834862385321101
$ws.Range("P2").Value = 0.1422018348623853
$ws.Range("S2").Value = 0.1284403669724771
# Row 3
$ws.Range("B3").Value = 0.00909090909090909
$ws.Range("J3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.6636363636363637
$ws.Range("S3").Value = 0.2818181818181818
# Row 4
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.7111111111111111
$ws.Range("S4").Value = 0.2666666666666667
# Row 6
$ws.Range("B6").Value = 0.03686635944700461
$ws.Range("D6").Value = 0.0184331797235023
$ws.Range("F6").Value = 0.03225806451612903
$ws.Range("J6").Value = 0.3824884792626728
$ws.Range("O6").Value = 0.009216589861751152
$ws.Range("Q6").Value = 0.1935483870967742
$ws.Range("R6").Value = 0.05990783410138249
$ws.Range("S6").Value = 0.2672811059907834
# Row 7
$ws.Range("B7").Value = 0.131578947368421
$ws.Range("D7").Value = 0.0131578947368421
$ws.Range("F7").Value = 0.05921052631578947
$ws.Range("J7").Value = 0.1052631578947368
$ws.Range("O7").Value = 0.006578947368421052
$ws.Range("Q7").Value = 0.1776315789473684
$ws.Range("R7").Value = 0.1052631578947368
$ws.Range("S7").Value = 0.4013157894736842
# Row 8
$ws.Range("B8").Value = 0.0585480093676815
$ws.Range("D8").Value = 0.0234192037470726
$ws.Range("F8").Value = 0.06557377049180328
$ws.Range("J8").Value = 0.117096018735363
$ws.Range("O8").Value = 0.0117096018735363
$ws.Range("Q8").Value = 0.1756440281030445
$ws.Range("R8").Value = 0.1545667447306792
$ws.Range("S8").Value = 0.3934426229508197
# Row 9
$ws.Range("B9").Value = 0.05633802816901409
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.07042253521126761
$ws.Range("J9").Value = 0.09859154929577464
$ws.Range("O9").Value = 0.02816901408450704
$ws.Range("Q9").Value = 0.2112676056338028
$ws.Range("R9").Value = 0.1549295774647887
$ws.Range("S9").Value = 0.3661971830985916
# Row 10
$ws.Range("B10").Value = 0.09067579127459367
$ws.Range("D10").Value = 0.02309666381522669
$ws.Range("E10").Value = 0.000855431993156544
$ws.Range("F10").Value = 0.0641573994867408
$ws.Range("J10").Value = 0.1120615911035073
$ws.Range("O10").Value = 0.01454234388366125
$ws.Range("Q10").Value = 0.1907613344739093
$ws.Range("R10").Value = 0.1240376390076989
$ws.Range("S10").Value = 0.3798118049615056
# Row 13
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.3055555555555556
# Row 15
$ws.Range("F15").Value = 0.0297029702970297
$ws.Range("H15").Value = 0.1881188118811881
$ws.Range("I15").Value = 0.0396039603960396
$ws.Range("J15").Value = 0.4158415841584158
$ws.Range("K15").Value = 0.0594059405940594
$ws.Range("M15").Value = 0.0297029702970297
$ws.Range("O15").Value = 0.03465346534653466
$ws.Range("S15").Value = 0.202970297029703
# Row 16
$ws.Range("F16").Value = 0.03759398496240601
$ws.Range("H16").Value = 0.2406015037593985
$ws.Range("I16").Value = 0.08270676691729323
$ws.Range("J16").Value = 0.3233082706766917
$ws.Range("K16").Value = 0.1278195488721804
$ws.Range("M16").Value = 0.02255639097744361
$ws.Range("O16").Value = 0.03007518796992481
$ws.Range("S16").Value = 0.1353383458646616
# Row 17
$ws.Range("F17").Value = 0.0303030303030303
$ws.Range("H17").Value = 0.1919191919191919
$ws.Range("I17").Value = 0.05555555555555555
$ws.Range("J17").Value = 0.4545454545454545
$ws.Range("K17").Value = 0.1136363636363636
$ws.Range("M17").Value = 0.005050505050505051
$ws.Range("O17").Value = 0.06818181818181818
$ws.Range("S17").Value = 0.08080808080808081
# Row 18
$ws.Range("F18").Value = 0.05343511450381679
$ws.Range("H18").Value = 0.1564885496183206
$ws.Range("I18").Value = 0.04961832061068702
$ws.Range("J18").Value = 0.4923664122137404
$ws.Range("K18").Value = 0.0916030534351145
$ws.Range("M18").Value = 0.01145038167938931
$ws.Range("O18").Value = 0.0648854961832061
$ws.Range("S18").Value = 0.08015267175572519
# Row 19
$ws.Range("F19").Value = 0.03192488262910798
$ws.Range("H19").Value = 0.2291079812206573
$ws.Range("I19").Value = 0.08262910798122065
$ws.Range("J19").Value = 0.3643192488262911
$ws.Range("K19").Value = 0.08262910798122065
$ws.Range("M19").Value = 0.0215962441314554
$ws.Range("N19").Value = 0.001877934272300469
$ws.Range("O19").Value = 0.07981220657276995
$ws.Range("S19").Value = 0.1061032863849765
